# "Test data for Greece Market"
# Duplicate the last country sheet (Croatia) to create a new "Greece" sheet,
# positioned right after it (i.e. at the end of the workbook), then update
# its market name / part-number cells.

$wb = $excel.ActiveWorkbook
$croatia = $wb.Worksheets.Item("Croatia")

# Copy Croatia immediately after itself -> becomes the new last sheet and
# the active sheet/tab.
$croatia.Copy($null, $croatia) | Out-Null

$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"

# Market-specific values (mirrors how every other country sheet is built).
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3206"

# Croatia keeps its data but loses its previous single-cell selection.
$croatia.Cells.Select() | Out-Null

# Greece is the newly active sheet with its own selection.
$greece.Select() | Out-Null
$greece.Range("D13").Select() | Out-Null
